# Insert a new data row at row 195 (pushing existing rows 195-283 down to
# 196-284) and populate the new row with the additional weekly price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(195).Insert()

$ws.Cells.Item(195, 1).Value  = 4
$ws.Cells.Item(195, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(195, 3).Value  = "Los Lagos"
$ws.Cells.Item(195, 4).Value  = 44609
$ws.Cells.Item(195, 5).Value  = 10
$ws.Cells.Item(195, 6).Value  = 100112023
$ws.Cells.Item(195, 7).Value  = "Brócoli"
$ws.Cells.Item(195, 8).Value  = "Sin especificar"
$ws.Cells.Item(195, 9).Value  = "Segunda"
$ws.Cells.Item(195, 10).Value = 250
$ws.Cells.Item(195, 11).Value = 1500
$ws.Cells.Item(195, 12).Value = 1500
$ws.Cells.Item(195, 13).Value = 1500
$ws.Cells.Item(195, 14).Value = "$/unidad"
$ws.Cells.Item(195, 15).Value = "Región Metropolitana"
$ws.Cells.Item(195, 16).Value = 1500
$ws.Cells.Item(195, 17).Value = 1
$ws.Cells.Item(195, 18).Value = "Hortaliza"
